$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.284.34"
$ws.Range("E2").Value = "  +0.20%  "

$ws.Range("D3").Value = "2.591.16"
$ws.Range("E3").Value = "  +1.66%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "588.54"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.53%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "148.98"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.39%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.02%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.597"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.93%  "

$ws.Range("E9").Value = "  +4.20%  "

$ws.Range("E10").Value = "  +1.81%  "

$ws.Range("E11").Value = "  +0.05%  "

$ws.Range("E12").Value = "  +1.60%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "27.75"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.82%  "

$ws.Range("D14").Value = "3.056.45"
$ws.Range("E14").Value = "  +1.80%  "

$ws.Range("D15").Value = "63.255.52"
$ws.Range("E15").Value = "  +0.32%  "

$ws.Range("E16").Value = "  +3.78%  "

$ws.Range("D17").Value = "2.573.45"
$ws.Range("E17").Value = "  +0.93%  "

$ws.Range("E18").Value = "  +0.25%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "345.17"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.91%  "

$ws.Range("E20").Value = "  +3.02%  "

$ws.Range("E21").Value = "  +1.77%  "

$ws.Range("E23").Value = "  -3.62%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "66.81"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.41%  "

$ws.Range("E25").Value = "  +0.52%  "

$ws.Range("D26").Value = "2.669.48"
$ws.Range("E26").Value = "  -0.13%  "

$ws.Range("E27").Value = "  -0.46%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.27"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +12.35%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.54"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.04%  "

$ws.Range("E30").Value = "  +0.13%  "

$ws.Range("E31").Value = "  -0.28%  "

$ws.Range("E32").Value = "  +7.56%  "

$ws.Range("E33").Value = "  +1.43%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "468.17"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +13.27%  "

$ws.Range("E35").Value = "  +4.58%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "176.94"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.48%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.407"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.82%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "19.31"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.92%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.67"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +6.31%  "

$ws.Range("E40").Value = "  +0.03%  "

$ws.Range("E41").Value = "  +0.46%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "152.39"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.21%  "

$ws.Range("E44").Value = "  +2.36%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "21.15"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.49%  "

$ws.Range("E46").Value = "  +5.58%  "

$ws.Range("E47").Value = "  +2.08%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0978"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.56%  "

$ws.Range("E49").Value = "  +1.44%  "

$ws.Range("E50").Value = "  -0.63%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "11.40"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.66%  "

